$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5564258098602295
$ws.Range("B1").Value = 3.923957824707031
$ws.Range("C1").Value = 6.129853248596191
$ws.Range("D1").Value = 1.459057688713074
$ws.Range("E1").Value = 0.8451275825500488
